# Leave 3/9/2023 12:08 AM
# Adds a new table row (table grows from A8:K130 to A8:K131), inserts a
# "2023" year-marker row, fills in monthly PERIOD dates for the new rows,
# and records a couple of leave entries (EARNED / SICK LEAVE).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# 1) Grow the table by one row (A8:K130 -> A8:K131). This keeps the table
#    definition / calculated columns consistent with the new row count.
$lo.ListRows.Add() | Out-Null

# 2) The table's special "last row" formatting (styles used only on the very
#    last row of the table) needs to move from row 130 down to the new row
#    131, and row 130 goes back to the regular body-row formatting (copied
#    from row 129, an ordinary interior row).
$ws.Range("A130:K130").Copy()
$ws.Range("A131:K131").PasteSpecial(-4122) | Out-Null
$ws.Range("A129:K129").Copy()
$ws.Range("A130:K130").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# The calculated column's formula doesn't automatically propagate onto the
# freshly-materialized row 131 from the format-only paste above, so restate
# it explicitly (matches every other body row's EARNED-mirror formula).
$ws.Range("G131").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# 3) Row 91 becomes the "2023" year-divider row, matching the style used for
#    the previous year dividers (e.g. row 77's "2022"). Copy that cell's
#    format, then force the text "2023" (not the number 2023) via a leading
#    apostrophe, the same way Excel stores these quote-prefixed year labels.
$ws.Range("A77").Copy()
$ws.Range("A91").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A91").Value = "'2023"

# 4) Fill in the monthly PERIOD date for each row from Feb-2023 through
#    Dec-2025 (rows 92-127), plus the two leave entries recorded in
#    Feb/Mar 2023 (rows 92 and 93).
$ws.Range("A92").Value = 44927
$ws.Range("C92").Value = 1.25

$ws.Range("A93").Value = 44958
$ws.Range("B93").Value = "SL(1-0-00)"
$ws.Range("C93").Value = 1.25
$ws.Range("H93").Value = 1

# K93 carries a date value with the same date-formatted style used by the
# other "as of" dates in column K (e.g. K85); copy that format over first.
$ws.Range("K85").Copy()
$ws.Range("K93").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("K93").Value = 44967

$ws.Range("A94").Value = 44986
$ws.Range("A95").Value = 45017
$ws.Range("A96").Value = 45047
$ws.Range("A97").Value = 45078
$ws.Range("A98").Value = 45108
$ws.Range("A99").Value = 45139
$ws.Range("A100").Value = 45170
$ws.Range("A101").Value = 45200
$ws.Range("A102").Value = 45231
$ws.Range("A103").Value = 45261
$ws.Range("A104").Value = 45292
$ws.Range("A105").Value = 45323
$ws.Range("A106").Value = 45352
$ws.Range("A107").Value = 45383
$ws.Range("A108").Value = 45413
$ws.Range("A109").Value = 45444
$ws.Range("A110").Value = 45474
$ws.Range("A111").Value = 45505
$ws.Range("A112").Value = 45536
$ws.Range("A113").Value = 45566
$ws.Range("A114").Value = 45597
$ws.Range("A115").Value = 45627
$ws.Range("A116").Value = 45658
$ws.Range("A117").Value = 45689
$ws.Range("A118").Value = 45717
$ws.Range("A119").Value = 45748
$ws.Range("A120").Value = 45778
$ws.Range("A121").Value = 45809
$ws.Range("A122").Value = 45839
$ws.Range("A123").Value = 45870
$ws.Range("A124").Value = 45901
$ws.Range("A125").Value = 45931
$ws.Range("A126").Value = 45962
$ws.Range("A127").Value = 45992

# 5) Reflect the last on-screen selection from the edit session.
$ws.Range("B94").Select()
